$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B16").Value = 'Pontuação: 1; Força = "Easy"'
$ws.Range("B16").Select()
